$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C takes over the old "Acceleration_SMA" header (with its bold/border/centered style).
# Column B becomes the new "Trening" (training type) column.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C1").Value = "Acceleration_SMA"
$ws.Range("B1").Value = "Trening"

# Row 2 (originally "10-15" / 3.925750268684639) is split into two training-type rows.
$ws.Range("A2").Value = "10-15"
$ws.Range("B2").Value = "Duża Gra"
$ws.Range("C2").Value = 3.841688387559884

$ws.Range("A3").Value = "10-15"
$ws.Range("B3").Value = "Mała Gra"
$ws.Range("C3").Value = 3.389416957949544

# Row 3 (originally "5-10" / 3.350359798132719) is split into two training-type rows.
$ws.Range("A4").Value = "5-10"
$ws.Range("B4").Value = "Duża Gra"
$ws.Range("C4").Value = 3.266047984689147

$ws.Range("A5").Value = "5-10"
$ws.Range("B5").Value = "Mała Gra"
$ws.Range("C5").Value = 2.928736258011598
